$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the contents of C1:H1 one column to the left (into B1:G1), effectively
# deleting the PHONE header and shifting everything after it left.
for ($col = 3; $col -le 8; $col++) {
    $srcCell = $ws.Cells.Item(1, $col)
    $dstCell = $ws.Cells.Item(1, $col - 1)
    $dstCell.Value = $srcCell.Value2
}

# Clear the now-trailing (duplicate) last cell's content, but keep its style
$ws.Cells.Item(1, 8).ClearContents()

# Update the active selection to match the post-edit state (G1)
$ws.Range("G1").Select()
